# Refined metadata to be additional tab
#
# 1) Refresh the "time_taken" timestamps in column F of the "data" sheet
#    (rows 2-48) to reflect a re-run of the panel query.
# 2) Add a new "metadata" worksheet (after "data") describing the panel
#    query that produced the data: data_name, data_id, data_version,
#    data_version_created, panel_query_time, panel_get_request.

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item(1)

# --- 1) Update the per-row query timestamps on the "data" sheet -----------
$timeTaken = @{
    2  = "2021-10-05 14:19:38.017089"
    3  = "2021-10-05 14:19:38.017098"
    4  = "2021-10-05 14:19:38.017102"
    5  = "2021-10-05 14:19:38.017105"
    6  = "2021-10-05 14:19:38.017108"
    7  = "2021-10-05 14:19:38.017111"
    8  = "2021-10-05 14:19:38.017113"
    9  = "2021-10-05 14:19:38.017116"
    10 = "2021-10-05 14:19:38.017119"
    11 = "2021-10-05 14:19:38.017122"
    12 = "2021-10-05 14:19:38.017124"
    13 = "2021-10-05 14:19:38.017127"
    14 = "2021-10-05 14:19:38.017130"
    15 = "2021-10-05 14:19:38.017132"
    16 = "2021-10-05 14:19:38.017135"
    17 = "2021-10-05 14:19:38.017138"
    18 = "2021-10-05 14:19:38.017141"
    19 = "2021-10-05 14:19:38.017144"
    20 = "2021-10-05 14:19:38.017147"
    21 = "2021-10-05 14:19:38.017150"
    22 = "2021-10-05 14:19:38.017152"
    23 = "2021-10-05 14:19:38.017155"
    24 = "2021-10-05 14:19:38.017157"
    25 = "2021-10-05 14:19:38.017160"
    26 = "2021-10-05 14:19:38.017163"
    27 = "2021-10-05 14:19:38.017166"
    28 = "2021-10-05 14:19:38.017169"
    29 = "2021-10-05 14:19:38.017171"
    30 = "2021-10-05 14:19:38.017174"
    31 = "2021-10-05 14:19:38.017176"
    32 = "2021-10-05 14:19:38.017179"
    33 = "2021-10-05 14:19:38.017182"
    34 = "2021-10-05 14:19:38.017185"
    35 = "2021-10-05 14:19:38.017187"
    36 = "2021-10-05 14:19:38.017190"
    37 = "2021-10-05 14:19:38.017193"
    38 = "2021-10-05 14:19:38.017196"
    39 = "2021-10-05 14:19:38.017198"
    40 = "2021-10-05 14:19:38.017201"
    41 = "2021-10-05 14:19:38.017204"
    42 = "2021-10-05 14:19:38.017207"
    43 = "2021-10-05 14:19:38.017210"
    44 = "2021-10-05 14:19:38.017213"
    45 = "2021-10-05 14:19:38.017215"
    46 = "2021-10-05 14:19:38.017218"
    47 = "2021-10-05 14:19:38.017221"
    48 = "2021-10-05 14:19:38.017223"
}

foreach ($row in $timeTaken.Keys) {
    $dataSheet.Cells.Item($row, 6).Value = $timeTaken[$row]
}

# --- 2) Add the "metadata" sheet, placed right after "data" ---------------
$metaSheet = $wb.Worksheets.Add($null, $dataSheet)
$metaSheet.Name = "metadata"

# Header row (bold, centered, top-aligned, thin border - same look as the
# "data" sheet's header row).
$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

$header = $metaSheet.Range("B1:G1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160
$header.Borders.Item(7).LineStyle = 1
$header.Borders.Item(8).LineStyle = 1
$header.Borders.Item(9).LineStyle = 1
$header.Borders.Item(10).LineStyle = 1

# Data row.
$metaSheet.Range("A2").Value = 0
$metaSheet.Range("A2").Font.Bold = $true
$metaSheet.Range("A2").HorizontalAlignment = -4108
$metaSheet.Range("A2").VerticalAlignment = -4160
$metaSheet.Range("A2").Borders.Item(7).LineStyle = 1
$metaSheet.Range("A2").Borders.Item(8).LineStyle = 1
$metaSheet.Range("A2").Borders.Item(9).LineStyle = 1
$metaSheet.Range("A2").Borders.Item(10).LineStyle = 1

$metaSheet.Range("B2").Value = "ClinGen Gene Validity Curations"
$metaSheet.Range("C2").Value = 64

$versionCell = $metaSheet.Range("D2")
$versionCell.NumberFormat = "@"
$versionCell.Value = "0.64"

$metaSheet.Range("E2").Value = "2019-06-20T15:10:34.572009Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:19:38.013619"
$metaSheet.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/64/?format=json"

# Keep "data" as the active sheet/tab, same as before the edit.
$dataSheet.Activate()
